$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1525498355650114
$ws.Range("D2").Value = 0.01978436094723257
$ws.Range("E2").Value = 0.09120834119893473
$ws.Range("F2").Value = 1.039537259670652
$ws.Range("G2").Value = 0.002425020479432729
$ws.Range("I2").Value = 0.9731179099746328
$ws.Range("L2").Value = 0.1082803219544983
$ws.Range("M2").Value = 1.196563620267767
$ws.Range("N2").Value = 1.774053148430852
$ws.Range("O2").Value = 3.624946256139026
$ws.Range("C3").Value = 0.154260758059241
$ws.Range("D3").Value = 0.02000544761482859
$ws.Range("E3").Value = 0.09286899463377518
$ws.Range("F3").Value = 1.011003328889217
$ws.Range("G3").Value = 0.002429262197577043
$ws.Range("I3").Value = 0.9563343663396182
$ws.Range("L3").Value = 0.1105012220155501
$ws.Range("M3").Value = 1.079975144303546
$ws.Range("N3").Value = 1.633127766037092
$ws.Range("O3").Value = 3.541073715263963
$ws.Range("C4").Value = 0.1553830669044487
$ws.Range("D4").Value = 0.02014725347299073
$ws.Range("E4").Value = 0.09394517187478457
$ws.Range("F4").Value = 0.9942231773205066
$ws.Range("G4").Value = 0.002432006603751962
$ws.Range("I4").Value = 0.9467221436595921
$ws.Range("L4").Value = 0.1119419180728993
$ws.Range("M4").Value = 1.008250412120361
$ws.Range("N4").Value = 1.54679984475311
$ws.Range("O4").Value = 3.492164271692445
$ws.Range("C5").Value = 0.1558584556828642
$ws.Range("D5").Value = 0.02020656256090136
$ws.Range("E5").Value = 0.09439791282622079
$ws.Range("F5").Value = 0.9875704221649499
$ws.Range("G5").Value = 0.002433160283143916
$ws.Range("I5").Value = 0.9429786942008249
$ws.Range("L5").Value = 0.1125483527025919
$ws.Range("M5").Value = 0.9789896064578585
$ws.Range("N5").Value = 1.511675152813581
$ws.Range("O5").Value = 3.472881645203643
$ws.Range("C6").Value = 0.1559384826967545
$ws.Range("D6").Value = 0.02021650268407882
$ws.Range("E6").Value = 0.09447394643256513
$ws.Range("F6").Value = 0.9864769070927082
$ws.Range("G6").Value = 0.002433353986788789
$ws.Range("I6").Value = 0.9423675643306311
$ws.Range("L6").Value = 0.1126502176869106
$ws.Range("M6").Value = 0.9741289904716126
$ws.Range("N6").Value = 1.505846146994145
$ws.Range("O6").Value = 3.469718865053323
$ws.Range("C7").Value = 0.1553894051448239
$ws.Range("D7").Value = 0.0201480471747888
$ws.Range("E7").Value = 0.09395122027622116
$ws.Range("F7").Value = 0.994132706671266
$ws.Range("G7").Value = 0.00243202201952586
$ws.Range("I7").Value = 0.9466709561000286
$ws.Range("L7").Value = 0.111950018407768
$ws.Range("M7").Value = 1.00785591845117
$ws.Range("N7").Value = 1.546325914133831
$ws.Range("O7").Value = 3.491901597638559
$ws.Range("C8").Value = 0.1531248517285348
$ws.Range("D8").Value = 0.01985933360777503
$ws.Range("E8").Value = 0.09176918714206983
$ws.Range("F8").Value = 1.029544760222393
$ws.Range("G8").Value = 0.00242645404098496
$ws.Range("I8").Value = 0.9671867323828565
$ws.Range("L8").Value = 0.109030073799568
$ws.Range("M8").Value = 1.156394234190017
$ws.Range("N8").Value = 1.725422813166972
$ws.Range("O8").Value = 3.595488149417974
$ws.Range("C9").Value = 0.1492543595368403
$ws.Range("D9").Value = 0.01934128095765519
$ws.Range("E9").Value = 0.08793979938449958
$ws.Range("F9").Value = 1.104895912507928
$ws.Range("G9").Value = 0.002416640683719813
$ws.Range("I9").Value = 1.012948539358149
$ws.Range("L9").Value = 0.1039170208396811
$ws.Range("M9").Value = 1.446477761939889
$ws.Range("N9").Value = 2.078062853888412
$ws.Range("O9").Value = 3.819292285386211
$ws.Range("C10").Value = 0.1467591805051605
$ws.Range("D10").Value = 0.01899005629909123
$ws.Range("E10").Value = 0.08540174195976036
$ws.Range("F10").Value = 1.163918541516125
$ws.Range("G10").Value = 0.002410097434924348
$ws.Range("I10").Value = 1.04999101886601
$ws.Range("L10").Value = 0.1005361466303487
$ws.Range("M10").Value = 1.658761291822785
$ws.Range("N10").Value = 2.337831498260812
$ws.Range("O10").Value = 3.996527943798753
$ws.Range("C11").Value = 0.1456999015387801
$ws.Range("D11").Value = 0.01883667085239171
$ws.Range("E11").Value = 0.08430719513995166
$ws.Range("F11").Value = 1.191578176916337
$ws.Range("G11").Value = 0.002407263942888237
$ws.Range("I11").Value = 1.067596603684876
$ws.Range("L11").Value = 0.09908010975093795
$ws.Range("M11").Value = 1.755130312702761
$ws.Range("N11").Value = 2.456118070216178
$ws.Range("O11").Value = 4.079983729685296
$ws.Range("C12").Value = 0.1453096959324895
$ws.Range("D12").Value = 0.01877950833886644
$ws.Range("E12").Value = 0.08390137754951044
$ws.Range("F12").Value = 1.202169581104357
$ws.Range("G12").Value = 0.002406211427120591
$ws.Range("I12").Value = 1.07437272063693
$ws.Range("L12").Value = 0.09854056995914462
$ws.Range("M12").Value = 1.791591770605237
$ws.Range("N12").Value = 2.500923245998706
$ws.Range("O12").Value = 4.111996547542901
$ws.Range("C13").Value = 0.1453932476579283
$ws.Range("D13").Value = 0.01879177829509704
$ws.Range("E13").Value = 0.08398839181126849
$ws.Range("F13").Value = 1.199883302424908
$ws.Range("G13").Value = 0.002406437196860676
$ws.Range("I13").Value = 1.072908491403581
$ws.Range("L13").Value = 0.09865624262289607
$ws.Range("M13").Value = 1.783740581943846
$ws.Range("N13").Value = 2.491273153794225
$ws.Range("O13").Value = 4.105083736321717
$ws.Range("C14").Value = 0.1456675801029235
$ws.Range("D14").Value = 0.01883194957971224
$ws.Range("E14").Value = 0.0842736344711692
$ws.Range("F14").Value = 1.192447182450806
$ws.Range("G14").Value = 0.00240717694227616
$ws.Range("I14").Value = 1.068151884898711
$ws.Range("L14").Value = 0.09903548413136498
$ws.Range("M14").Value = 1.758130663257219
$ws.Range("N14").Value = 2.459803989600573
$ws.Range("O14").Value = 4.082609211229624
$ws.Range("C15").Value = 0.1458370394483168
$ws.Range("D15").Value = 0.01885667570268712
$ws.Range("E15").Value = 0.08444948302381994
$ws.Range("F15").Value = 1.187907645031132
$ws.Range("G15").Value = 0.002407632719760968
$ws.Range("I15").Value = 1.06525257718738
$ws.Range("L15").Value = 0.09926932249090825
$ws.Range("M15").Value = 1.742439680806712
$ws.Range("N15").Value = 2.440529739330714
$ws.Range("O15").Value = 4.068896386660981
$ws.Range("C16").Value = 0.1468299338323114
$ws.Range("D16").Value = 0.01900020916822065
$ws.Range("E16").Value = 0.08547448378292621
$ws.Range("F16").Value = 1.162127288389755
$ws.Range("G16").Value = 0.002410285480126641
$ws.Range("I16").Value = 1.048855702632622
$ws.Range("L16").Value = 0.1006329549054676
$ws.Range("M16").Value = 1.652459100482218
$ws.Range("N16").Value = 2.330103182268999
$ws.Range("O16").Value = 3.991131115424707
$ws.Range("C17").Value = 0.1474584695776748
$ws.Range("D17").Value = 0.01908990031404478
$ws.Range("E17").Value = 0.08611868561126734
$ws.Range("F17").Value = 1.146519910599409
$ws.Range("G17").Value = 0.002411949429470463
$ws.Range("I17").Value = 1.038990583677659
$ws.Range("L17").Value = 0.1014905192327209
$ws.Range("M17").Value = 1.597205847203341
$ws.Range("N17").Value = 2.262386947027665
$ws.Range("O17").Value = 3.944151639810627
$ws.Range("C18").Value = 0.1478271191632281
$ws.Range("D18").Value = 0.01914208943375773
$ws.Range("E18").Value = 0.08649486171085652
$ws.Range("F18").Value = 1.137619145516481
$ws.Range("G18").Value = 0.002412919961437515
$ws.Range("I18").Value = 1.033387427321671
$ws.Range("L18").Value = 0.1019914772818709
$ws.Range("M18").Value = 1.565406987116177
$ws.Range("N18").Value = 2.223449569927141
$ws.Range("O18").Value = 3.917396522627882
$ws.Range("C19").Value = 0.1479531619228283
$ws.Range("D19").Value = 0.01915986294317751
$ws.Range("E19").Value = 0.08662319760708737
$ws.Range("F19").Value = 1.134618559148151
$ws.Range("G19").Value = 0.002413250883862403
$ws.Range("I19").Value = 1.031502466916265
$ws.Range("L19").Value = 0.1021624157981567
$ws.Range("M19").Value = 1.554637317822397
$ws.Range("N19").Value = 2.210268095428091
$ws.Range("O19").Value = 3.908383337230532
$ws.Range("C20").Value = 0.1473908225361065
$ws.Range("D20").Value = 0.01908029030455882
$ws.Range("E20").Value = 0.08604952442850622
$ws.Range("F20").Value = 1.148173451680861
$ws.Range("G20").Value = 0.002411770905746991
$ws.Range("I20").Value = 1.040033389048446
$ws.Range("L20").Value = 0.1013984317823375
$ws.Range("M20").Value = 1.603089599930357
$ws.Range("N20").Value = 2.269594325780304
$ws.Range("O20").Value = 3.949125109935835
$ws.Range("C21").Value = 0.1455867054646198
$ws.Range("D21").Value = 0.01882012527689003
$ws.Range("E21").Value = 0.08418961640940148
$ws.Range("F21").Value = 1.194628160979917
$ws.Range("G21").Value = 0.002406959106335062
$ws.Range("I21").Value = 1.069546044303223
$ws.Range("L21").Value = 0.09892377027469124
$ws.Range("M21").Value = 1.765653790456867
$ws.Range("N21").Value = 2.469046936180973
$ws.Range("O21").Value = 4.089199375497287
$ws.Range("C22").Value = 0.1444712730345934
$ws.Range("D22").Value = 0.01865546330179235
$ws.Range("E22").Value = 0.08302457147862763
$ws.Range("F22").Value = 1.225673106805559
$ws.Range("G22").Value = 0.002403933563995979
$ws.Range("I22").Value = 1.089471487181044
$ws.Range("L22").Value = 0.09737541095816837
$ws.Range("M22").Value = 1.871714648326787
$ws.Range("N22").Value = 2.599472415249352
$ws.Range("O22").Value = 4.183136719767958
$ws.Range("C23").Value = 0.1450607679696603
$ws.Range("D23").Value = 0.01874285416769261
$ws.Range("E23").Value = 0.08364174552152137
$ws.Range("F23").Value = 1.209040962668809
$ws.Range("G23").Value = 0.002405537476485511
$ws.Range("I23").Value = 1.078778355874761
$ws.Range("L23").Value = 0.0981954726378067
$ws.Range("M23").Value = 1.8151257052576
$ws.Range("N23").Value = 2.52985667699329
$ws.Range("O23").Value = 4.132780871424416
$ws.Range("C24").Value = 0.1474213830393936
$ws.Range("D24").Value = 0.01908463304398644
$ws.Range("E24").Value = 0.08608077408564863
$ws.Range("F24").Value = 1.147425661161094
$ws.Range("G24").Value = 0.002411851572970414
$ws.Range("I24").Value = 1.039561723636595
$ws.Range("L24").Value = 0.1014400398118394
$ws.Range("M24").Value = 1.600429657997168
$ws.Range("N24").Value = 2.266335889499715
$ws.Range("O24").Value = 3.946875813098472
$ws.Range("C25").Value = 0.1502402939102794
$ws.Range("D25").Value = 0.0194762716353658
$ws.Range("E25").Value = 0.08892748440155385
$ws.Range("F25").Value = 1.083872700377412
$ws.Range("G25").Value = 0.002419177867402225
$ws.Range("I25").Value = 0.9999715613486302
$ws.Range("L25").Value = 0.1052344098808771
$ws.Range("M25").Value = 1.368142308728864
$ws.Range("N25").Value = 1.982530848130693
$ws.Range("O25").Value = 3.756513228528092
